$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: "_old" suffix -> "_FV2410", "_new" suffix -> "_FV2504"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace '_old$', '_FV2410'
        $newVal = $newVal -replace '_new$', '_FV2504'
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}

# Turn the data range A1:U76 into an Excel Table ("Table1") with an autofilter
$range = $ws.Range("A1:U76")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row (split below row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$wb.Save()
